# TC10_INS_CancerType-Leukemia.xlsx edit
# - Rewrite the "ProgramsTab" query (B2) so the "Website" column is now
#   computed via a CASE expression (program_link / program_acronym) instead
#   of a straight prg.website reference.
# - Update the sheet's saved cursor/selection from C3 to C8 (and best-effort
#   scroll the viewport so row 5 is back at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProgramsQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Leukemia%'
ORDER BY 
   lower(prg.program_name) ASC
LIMIT 100;
"@

# The here-string adds a trailing newline; the stored cell text does not
# end with one, so trim it back off.
$newProgramsQuery = $newProgramsQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newProgramsQuery

# Restore the sheet's wrap-text / font-size formatting on the cell (keeps
# parity with the other query cells in column B).
$ws.Range("B2").WrapText = $true
$ws.Range("B2").Font.Size = 12

# Move the saved selection/cursor to C8 and scroll so row 5 is at the top,
# matching the workbook's last interactive state.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("C8").Select()
